$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '303.37'
Set-TextValue 'E2' '-1.02%'
Set-TextValue 'D4' '5.037'
Set-TextValue 'E4' '-1.36%'
Set-TextValue 'D5' '0.07891'
Set-TextValue 'E5' '-2.86%'
Set-TextValue 'D6' '1.827'
Set-TextValue 'E6' '-6.15%'
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue 'D7' '4.106'
Set-TextValue 'E7' '-2.14%'
$ws.Range('B8').Value = 'KuCoinToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
Set-TextValue 'D8' '7.785'
Set-TextValue 'E8' '0.08%'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D9' '0.9203'
Set-TextValue 'E9' '-1.38%'
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue 'D10' '0.1364'
Set-TextValue 'E10' '-2.27%'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue 'D11' '0.1895'
Set-TextValue 'E11' '-1.47%'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue 'D12' '0.09093'
Set-TextValue 'E12' '-1.59%'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue 'D13' '0.03469'
Set-TextValue 'E13' '-2.93%'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue 'D14' '0.09824'
Set-TextValue 'E14' '-0.31%'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue 'D15' '0.001422'
Set-TextValue 'E15' '0.85%'
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue 'D16' '0.006065'
Set-TextValue 'E16' '3.61%'
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 'D17' '3.722'
Set-TextValue 'E17' '3.75%'
Set-TextValue 'D18' '3.344'
Set-TextValue 'E18' '12.19%'
Set-TextValue 'E19' '-0.05%'
$ws.Range('B20').Value = 'MCDex'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue 'D20' '5.185'
Set-TextValue 'E20' '6.00%'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue 'D21' '0.1310'
Set-TextValue 'E21' '-2.92%'
Set-TextValue 'D23' '0.04414'
Set-TextValue 'E23' '-2.02%'
Set-TextValue 'D24' '0.001237'
Set-TextValue 'E24' '1.76%'
Set-TextValue 'D25' '0.004615'
Set-TextValue 'E25' '-5.31%'
Set-TextValue 'D26' '0.0001302'
Set-TextValue 'E26' '4.82%'
Set-TextValue 'D39' '0.01938'
Set-TextValue 'E39' '-3.11%'
Set-TextValue 'D40' '0.05076'
Set-TextValue 'E40' '2.86%'
Set-TextValue 'D41' '0.007623'
Set-TextValue 'E41' '-0.62%'
Set-TextValue 'D42' '0.01018'
Set-TextValue 'E42' '-8.60%'
Set-TextValue 'D43' '0.1343'
Set-TextValue 'E43' '-2.92%'
Set-TextValue 'D44' '0.002163'
Set-TextValue 'E44' '2.84%'
Set-TextValue 'E45' '-4.59%'
Set-TextValue 'D46' '0.00006198'
Set-TextValue 'E46' '-4.09%'
Set-TextValue 'E47' '-0.03%'
Set-TextValue 'D48' '65.22'
Set-TextValue 'E48' '0.85%'
Set-TextValue 'E49' '39.27%'
Set-TextValue 'D50' '0.00002103'
Set-TextValue 'E50' '-0.03%'
Set-TextValue 'D51' '0.0002003'
Set-TextValue 'E51' '-0.03%'
